# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 46061 (2026-02-08) to 46062 (2026-02-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 3).Value = 46062
}
